# Scheduled runner update: refresh cached Universalis market-price figures
# (currentAveragePrice*, Leve* profit columns) across the per-job Leve
# profit sheets. Values below are the newly fetched snapshot; columns are
# H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
# K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 104.333336
$ws.Cells.Item(11, 9).Value = 104.333336
$ws.Cells.Item(11, 11).Value = 104.333336
$ws.Cells.Item(11, 13).Value = 35.666664
$ws.Cells.Item(33, 8).Value = 185.04167
$ws.Cells.Item(33, 9).Value = 184.6
$ws.Cells.Item(33, 11).Value = 184.6
$ws.Cells.Item(33, 13).Value = 44.40000000000001
$ws.Cells.Item(53, 8).Value = 219.22223
$ws.Cells.Item(53, 9).Value = 218
$ws.Cells.Item(53, 11).Value = 218
$ws.Cells.Item(53, 13).Value = 419
$ws.Cells.Item(74, 8).Value = 12121.375
$ws.Cells.Item(74, 9).Value = 11567.286
$ws.Cells.Item(74, 11).Value = 11567.286
$ws.Cells.Item(74, 13).Value = -10631.286
$ws.Cells.Item(77, 8).Value = 12121.375
$ws.Cells.Item(77, 9).Value = 11567.286
$ws.Cells.Item(77, 11).Value = 57836.43
$ws.Cells.Item(77, 13).Value = -53156.43
$ws.Cells.Item(80, 8).Value = 1699.3334
$ws.Cells.Item(80, 10).Value = 1732.5
$ws.Cells.Item(80, 12).Value = 5197.5
$ws.Cells.Item(80, 14).Value = -7193.5
$ws.Cells.Item(83, 8).Value = 1699.3334
$ws.Cells.Item(83, 10).Value = 1732.5
$ws.Cells.Item(83, 12).Value = 15592.5
$ws.Cells.Item(83, 14).Value = -25576.5
$ws.Cells.Item(92, 8).Value = 1040.591
$ws.Cells.Item(92, 9).Value = 1004.2353
$ws.Cells.Item(92, 11).Value = 1004.2353
$ws.Cells.Item(92, 13).Value = 243.7646999999999
$ws.Cells.Item(96, 8).Value = 1287.75
$ws.Cells.Item(96, 9).Value = 345.4
$ws.Cells.Item(96, 10).Value = 5999.5
$ws.Cells.Item(96, 11).Value = 1036.2
$ws.Cells.Item(96, 12).Value = 17998.5
$ws.Cells.Item(96, 13).Value = 336.8000000000002
$ws.Cells.Item(96, 14).Value = -20744.5
$ws.Cells.Item(113, 8).Value = 6492.5
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 6492.5
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 6492.5
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -13000.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 21509.8
$ws.Cells.Item(23, 9).Value = 17666.334
$ws.Cells.Item(23, 11).Value = 17666.334
$ws.Cells.Item(23, 13).Value = -17407.334
$ws.Cells.Item(32, 8).Value = 19143
$ws.Cells.Item(32, 9).Value = 17819.8
$ws.Cells.Item(32, 10).Value = 32375
$ws.Cells.Item(32, 11).Value = 17819.8
$ws.Cells.Item(32, 12).Value = 32375
$ws.Cells.Item(32, 13).Value = -17532.8
$ws.Cells.Item(32, 14).Value = -32949
$ws.Cells.Item(74, 8).Value = 5719.154
$ws.Cells.Item(74, 9).Value = 7650.222
$ws.Cells.Item(74, 11).Value = 7650.222
$ws.Cells.Item(74, 13).Value = -6776.222
$ws.Cells.Item(77, 8).Value = 5719.154
$ws.Cells.Item(77, 9).Value = 7650.222
$ws.Cells.Item(77, 11).Value = 38251.11
$ws.Cells.Item(77, 13).Value = -33883.11
$ws.Cells.Item(131, 8).Value = 106666.336
$ws.Cells.Item(131, 10).Value = 106666.336
$ws.Cells.Item(131, 12).Value = 106666.336
$ws.Cells.Item(131, 14).Value = -116746.336
$ws.Cells.Item(132, 8).Value = 4465.8335
$ws.Cells.Item(132, 9).Value = 3932.6667
$ws.Cells.Item(132, 11).Value = 11798.0001
$ws.Cells.Item(132, 13).Value = -9268.000100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(9, 8).Value = 90000
$ws.Cells.Item(9, 10).Value = 90000
$ws.Cells.Item(9, 12).Value = 90000
$ws.Cells.Item(9, 14).Value = -90336
$ws.Cells.Item(105, 8).Value = 1881
$ws.Cells.Item(105, 9).Value = 1881
$ws.Cells.Item(105, 11).Value = 1881
$ws.Cells.Item(105, 13).Value = -134
$ws.Cells.Item(127, 8).Value = 49750
$ws.Cells.Item(127, 10).Value = 49750
$ws.Cells.Item(127, 12).Value = 49750
$ws.Cells.Item(127, 14).Value = -59670
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 40078.5
$ws.Cells.Item(50, 10).Value = 40078.5
$ws.Cells.Item(50, 12).Value = 40078.5
$ws.Cells.Item(50, 14).Value = -41328.5
$ws.Cells.Item(60, 8).Value = 35642
$ws.Cells.Item(60, 10).Value = 48014
$ws.Cells.Item(60, 12).Value = 48014
$ws.Cells.Item(60, 14).Value = -49036
$ws.Cells.Item(74, 8).Value = 69709.336
$ws.Cells.Item(74, 10).Value = 69709.336
$ws.Cells.Item(74, 12).Value = 69709.336
$ws.Cells.Item(74, 14).Value = -71457.336
$ws.Cells.Item(77, 8).Value = 69709.336
$ws.Cells.Item(77, 10).Value = 69709.336
$ws.Cells.Item(77, 12).Value = 209128.008
$ws.Cells.Item(77, 14).Value = -217864.008
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 4636.9165
$ws.Cells.Item(109, 9).Value = 2822
$ws.Cells.Item(109, 10).Value = 4999.9
$ws.Cells.Item(109, 11).Value = 8466
$ws.Cells.Item(109, 12).Value = 14999.7
$ws.Cells.Item(109, 13).Value = -7426
$ws.Cells.Item(109, 14).Value = -17079.7
$ws.Cells.Item(131, 8).Value = 3605.625
$ws.Cells.Item(131, 10).Value = 3611.1843
$ws.Cells.Item(131, 12).Value = 10833.5529
$ws.Cells.Item(131, 14).Value = -20913.5529
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6299.5
$ws.Cells.Item(22, 9).Value = 4432.3335
$ws.Cells.Item(22, 10).Value = 7699.875
$ws.Cells.Item(22, 11).Value = 4432.3335
$ws.Cells.Item(22, 12).Value = 7699.875
$ws.Cells.Item(22, 13).Value = -4137.3335
$ws.Cells.Item(22, 14).Value = -8289.875
$ws.Cells.Item(27, 8).Value = 6299.5
$ws.Cells.Item(27, 9).Value = 4432.3335
$ws.Cells.Item(27, 10).Value = 7699.875
$ws.Cells.Item(27, 11).Value = 4432.3335
$ws.Cells.Item(27, 12).Value = 7699.875
$ws.Cells.Item(27, 13).Value = -4325.3335
$ws.Cells.Item(55, 8).Value = 973.1667
$ws.Cells.Item(55, 9).Value = 596.6
$ws.Cells.Item(55, 10).Value = 1242.1428
$ws.Cells.Item(55, 11).Value = 596.6
$ws.Cells.Item(55, 12).Value = 1242.1428
$ws.Cells.Item(55, 13).Value = -423.6
$ws.Cells.Item(55, 14).Value = -1588.1428
$ws.Cells.Item(61, 8).Value = 1250
$ws.Cells.Item(61, 9).Value = 1000
$ws.Cells.Item(61, 11).Value = 1000
$ws.Cells.Item(61, 13).Value = -798
$ws.Cells.Item(82, 8).Value = 2659.4285
$ws.Cells.Item(82, 9).Value = 2443.2
$ws.Cells.Item(82, 10).Value = 3200
$ws.Cells.Item(82, 11).Value = 2443.2
$ws.Cells.Item(82, 12).Value = 3200
$ws.Cells.Item(82, 13).Value = -2082.2
$ws.Cells.Item(82, 14).Value = -3922
$ws.Cells.Item(85, 8).Value = 2659.4285
$ws.Cells.Item(85, 9).Value = 2443.2
$ws.Cells.Item(85, 10).Value = 3200
$ws.Cells.Item(85, 11).Value = 2443.2
$ws.Cells.Item(85, 12).Value = 3200
$ws.Cells.Item(85, 13).Value = -1195.2
$ws.Cells.Item(85, 14).Value = -5696
$ws.Cells.Item(93, 8).Value = 2010.3077
$ws.Cells.Item(93, 9).Value = 1919.75
$ws.Cells.Item(93, 11).Value = 1919.75
$ws.Cells.Item(93, 13).Value = -671.75
$ws.Cells.Item(113, 8).Value = 1250
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 13).Value = 1170
$ws.Cells.Item(127, 8).Value = 30000
$ws.Cells.Item(127, 10).Value = 30000
$ws.Cells.Item(127, 12).Value = 30000
$ws.Cells.Item(127, 14).Value = -39920
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 34613.332
$ws.Cells.Item(54, 10).Value = 34613.332
$ws.Cells.Item(54, 12).Value = 34613.332
$ws.Cells.Item(54, 14).Value = -35653.332
$ws.Cells.Item(123, 8).Value = 20000
$ws.Cells.Item(123, 9).Value = 20000
$ws.Cells.Item(123, 11).Value = 20000
$ws.Cells.Item(123, 13).Value = -15100
